# Update the "Förändrad" (Changed) date column (C) for all data rows.
# The tracked "changed" date serial moves from 46060 (2026-02-07) to
# 46061 (2026-02-08) for every data row (rows 2-41) on the only sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$lastRow = $ws.UsedRange.Rows.Count

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    $val = $cell.Value2()
    if ($val -eq 46060) {
        $cell.Value = 46061
    }
}
